$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1) "Leverandører" sheet (2nd sheet): fill in missing Bransje (C) / Webside
#    category (D) info, and correct the "For dyrt?" -> "Sjekk CE-merke..."
#    header text for two rows.
#
#    The new cell values are written in the same order the author appears to
#    have typed them in, so freshly-introduced shared strings land at the
#    same table indices as the target workbook (121..125).
# ---------------------------------------------------------------------------
$wsLev = $wb.Worksheets.Item("Leverandører")

$wsLev.Range("D17").Value = "Elektro / Automasjon"
$wsLev.Range("B3").Value = "Sjekk CE-merke / om leverandør er til å stole på"
$wsLev.Range("B4").Value = "Sjekk CE-merke / om leverandør er til å stole på"
$wsLev.Range("D27").Value = "Elektro / Automasjon / Elektronikk"
$wsLev.Range("D32").Value = "Automasjon"
$wsLev.Range("D36").Value = "Verktøy"

$wsLev.Range("C5").Value = "OK"
$wsLev.Range("C6").Value = "OK"
$wsLev.Range("C14").Value = "OK"
$wsLev.Range("C15").Value = "OK"
$wsLev.Range("D18").Value = "Blandet"
$wsLev.Range("D20").Value = "Elektro"
$wsLev.Range("D23").Value = "Elektro"
$wsLev.Range("D24").Value = "Elektro"
$wsLev.Range("C25").Value = "OK"
$wsLev.Range("C26").Value = "OK"
$wsLev.Range("C28").Value = "OK"
$wsLev.Range("D30").Value = "Elektro"
$wsLev.Range("C31").Value = "OK"
$wsLev.Range("C33").Value = "OK"
$wsLev.Range("C34").Value = "OK"
$wsLev.Range("C35").Value = "OK"

# ---------------------------------------------------------------------------
# 2) Remove the empty "Bøker 2ELR" sheet (3rd sheet) entirely.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Bøker 2ELR").Delete()

# ---------------------------------------------------------------------------
# 3) Restore the view state on the first sheet (scrolled down a bit, but
#    selection still on F15) ...
# ---------------------------------------------------------------------------
$wsFirst = $wb.Worksheets.Item("ElevregnskapUtstyr1ELR_Høst2018")
$wsFirst.Activate()
$wsFirst.Range("F15").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1

# ... and make "Leverandører" the active/selected tab, with D37 selected
# (the sheet that used to be 3rd - now removed - was the active tab before).
$wsLev.Activate()
$wsLev.Range("D37").Select()
